# Apply updated TPM-derived values to the NATMI LR-pair output sheet.
# The workbook holds ligand/receptor pair stats (Nlgn2-Nrxn2) across
# sending/target cluster combinations (rows 2-10, columns E-T).
# This reflects re-running the pipeline on updated TPM values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Sending=ECs, Target=ECs
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.164555333333333
$ws.Range("H2").Value = 3.493666
$ws.Range("I2").Value = 0.05923394707027321
$ws.Range("J2").Value = 0.05923394707027322
$ws.Range("M2").Value = 0.007742333333333334
$ws.Range("N2").Value = 0.023227
$ws.Range("O2").Value = 0.001217676423630818
$ws.Range("P2").Value = 0.001217676423630818
$ws.Range("Q2").Value = 0.009016375575777778
$ws.Range("R2").Value = 0.081147380182
$ws.Range("S2").Value = 0.00007212778082606745
$ws.Range("T2").Value = 0.00007212778082606745

# Row 3: Sending=ECs, Target=FAPs
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.164555333333333
$ws.Range("H3").Value = 3.493666
$ws.Range("I3").Value = 0.05923394707027321
$ws.Range("J3").Value = 0.05923394707027322
$ws.Range("O3").Value = 0.9497929577862038
$ws.Range("P3").Value = 0.9497929577862039
$ws.Range("Q3").Value = 7.032812543988
$ws.Range("R3").Value = 63.295312895892
$ws.Range("S3").Value = 0.05625998578922623
$ws.Range("T3").Value = 0.05625998578922625

# Row 4: Sending=ECs, Target=MuSCs
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.164555333333333
$ws.Range("H4").Value = 3.493666
$ws.Range("I4").Value = 0.05923394707027321
$ws.Range("J4").Value = 0.05923394707027322
$ws.Range("M4").Value = 0.3114883333333334
$ws.Range("N4").Value = 0.9344650000000001
$ws.Range("O4").Value = 0.04898936579016542
$ws.Range("P4").Value = 0.04898936579016543
$ws.Range("Q4").Value = 0.3627453998544445
$ws.Range("R4").Value = 3.26470859869
$ws.Range("S4").Value = 0.002901833500220911
$ws.Range("T4").Value = 0.002901833500220912

# Row 5: Sending=FAPs, Target=ECs
$ws.Range("I5").Value = 0.5488831985632208
$ws.Range("J5").Value = 0.5488831985632209
$ws.Range("M5").Value = 0.007742333333333334
$ws.Range("N5").Value = 0.023227
$ws.Range("O5").Value = 0.001217676423630818
$ws.Range("P5").Value = 0.001217676423630818
$ws.Range("Q5").Value = 0.08354900036644447
$ws.Range("R5").Value = 0.7519410032980002
$ws.Range("S5").Value = 0.0006683621302175067
$ws.Range("T5").Value = 0.0006683621302175068

# Row 6: Sending=FAPs, Target=FAPs
$ws.Range("I6").Value = 0.5488831985632208
$ws.Range("J6").Value = 0.5488831985632209
$ws.Range("O6").Value = 0.9497929577862038
$ws.Range("P6").Value = 0.9497929577862039
$ws.Range("S6").Value = 0.5213253966425136
$ws.Range("T6").Value = 0.5213253966425138

# Row 7: Sending=FAPs, Target=MuSCs
$ws.Range("I7").Value = 0.5488831985632208
$ws.Range("J7").Value = 0.5488831985632209
$ws.Range("M7").Value = 0.3114883333333334
$ws.Range("N7").Value = 0.9344650000000001
$ws.Range("O7").Value = 0.04898936579016542
$ws.Range("P7").Value = 0.04898936579016543
$ws.Range("R7").Value = 30.25197182791001
$ws.Range("S7").Value = 0.02688943979048962
$ws.Range("T7").Value = 0.02688943979048963

# Row 8: Sending=MuSCs, Target=ECs
$ws.Range("I8").Value = 0.391882854366506
$ws.Range("J8").Value = 0.3918828543665061
$ws.Range("M8").Value = 0.007742333333333334
$ws.Range("N8").Value = 0.023227
$ws.Range("O8").Value = 0.001217676423630818
$ws.Range("P8").Value = 0.001217676423630818
$ws.Range("Q8").Value = 0.05965098007877778
$ws.Range("R8").Value = 0.536858820709
$ws.Range("S8").Value = 0.0004771865125872437
$ws.Range("T8").Value = 0.0004771865125872437

# Row 9: Sending=MuSCs, Target=FAPs
$ws.Range("I9").Value = 0.391882854366506
$ws.Range("J9").Value = 0.3918828543665061
$ws.Range("O9").Value = 0.9497929577862038
$ws.Range("P9").Value = 0.9497929577862039
$ws.Range("S9").Value = 0.3722075753544639
$ws.Range("T9").Value = 0.372207575354464

# Row 10: Sending=MuSCs, Target=MuSCs
$ws.Range("I10").Value = 0.391882854366506
$ws.Range("J10").Value = 0.3918828543665061
$ws.Range("M10").Value = 0.3114883333333334
$ws.Range("N10").Value = 0.9344650000000001
$ws.Range("O10").Value = 0.04898936579016542
$ws.Range("P10").Value = 0.04898936579016543
$ws.Range("Q10").Value = 2.399868820739445
$ws.Range("S10").Value = 0.01919809249945489
$ws.Range("T10").Value = 0.01919809249945489
